$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 460.75
$ws.Range("I28").Value = 460.75
$ws.Range("K28").Value = 460.75
$ws.Range("M28").Value = 24.25

$ws.Range("H33").Value = 470.54544
$ws.Range("I33").Value = 532.7778
$ws.Range("J33").Value = 190.5
$ws.Range("K33").Value = 532.7778
$ws.Range("L33").Value = 190.5
$ws.Range("M33").Value = -303.7778
$ws.Range("N33").Value = -648.5

$ws.Range("H53").Value = 325.75
$ws.Range("I53").Value = 358.66666
$ws.Range("J53").Value = 306
$ws.Range("K53").Value = 358.66666
$ws.Range("L53").Value = 306
$ws.Range("M53").Value = 278.33334
$ws.Range("N53").Value = -1580

$ws.Range("H55").Value = 420.5
$ws.Range("I55").Value = 420.14285
$ws.Range("J55").Value = 421.33334
$ws.Range("K55").Value = 420.14285
$ws.Range("L55").Value = 421.33334
$ws.Range("M55").Value = -206.14285
$ws.Range("N55").Value = -849.33334

$ws.Range("H87").Value = 72975
$ws.Range("J87").Value = 72975
$ws.Range("L87").Value = 72975
$ws.Range("N87").Value = -75471

$ws.Range("H90").Value = 72975
$ws.Range("J90").Value = 72975
$ws.Range("L90").Value = 218925
$ws.Range("N90").Value = -231405

$ws.Range("H97").Value = 1478.875
$ws.Range("J97").Value = 1575.8572
$ws.Range("L97").Value = 4727.571599999999
$ws.Range("N97").Value = -5719.571599999999

$ws.Range("H98").Value = 907.38464
$ws.Range("I98").Value = 936.5454999999999
$ws.Range("J98").Value = 747
$ws.Range("K98").Value = 936.5454999999999
$ws.Range("L98").Value = 747
$ws.Range("M98").Value = 561.4545000000001
$ws.Range("N98").Value = -3743

$ws.Range("H107").Value = 1597.1666
$ws.Range("I107").Value = 1573.4706
$ws.Range("K107").Value = 1573.4706
$ws.Range("M107").Value = 346.5293999999999

$ws.Range("H111").Value = 1242.125
$ws.Range("I111").Value = 1308.8334
$ws.Range("J111").Value = 1042
$ws.Range("K111").Value = 3926.5002
$ws.Range("L111").Value = 3126
$ws.Range("M111").Value = -859.5001999999999
$ws.Range("N111").Value = -9260

$ws.Range("H122").Value = 907.38464
$ws.Range("I122").Value = 936.5454999999999
$ws.Range("J122").Value = 747
$ws.Range("K122").Value = 2809.6365
$ws.Range("L122").Value = 2241
$ws.Range("M122").Value = -359.6364999999996
$ws.Range("N122").Value = -7141

$ws.Range("H138").Value = 5260.3667
$ws.Range("J138").Value = 6763.381
$ws.Range("L138").Value = 20290.143
$ws.Range("N138").Value = -30570.143

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1954.3334
$ws.Range("I2").Value = 1954.3334
$ws.Range("K2").Value = 1954.3334
$ws.Range("M2").Value = -1841.3334

$ws.Range("H32").Value = 10572.483
$ws.Range("I32").Value = 9577.482
$ws.Range("K32").Value = 9577.482
$ws.Range("M32").Value = -9290.482

$ws.Range("H61").Value = 1952.3334
$ws.Range("I61").Value = 2153
$ws.Range("J61").Value = 1250
$ws.Range("K61").Value = 2153
$ws.Range("L61").Value = 1250
$ws.Range("M61").Value = -1941
$ws.Range("N61").Value = -1674

$ws.Range("H110").Value = 1876.0667
$ws.Range("I110").Value = 1393.7273
$ws.Range("J110").Value = 3202.5
$ws.Range("K110").Value = 1393.7273
$ws.Range("L110").Value = 3202.5
$ws.Range("M110").Value = 651.2727
$ws.Range("N110").Value = -7292.5

$ws.Range("H116").Value = 1954.3334
$ws.Range("I116").Value = 1954.3334
$ws.Range("K116").Value = 1954.3334
$ws.Range("M116").Value = 339.6666

$ws.Range("H132").Value = 3284.0667
$ws.Range("I132").Value = 2736.6
$ws.Range("J132").Value = 4379
$ws.Range("K132").Value = 8209.799999999999
$ws.Range("L132").Value = 13137
$ws.Range("M132").Value = -5679.799999999999
$ws.Range("N132").Value = -18197

$ws.Range("H136").Value = 1952.3334
$ws.Range("I136").Value = 2153
$ws.Range("J136").Value = 1250
$ws.Range("K136").Value = 6459
$ws.Range("L136").Value = 3750
$ws.Range("M136").Value = -3909
$ws.Range("N136").Value = -8850

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1954.3334
$ws.Range("I3").Value = 1954.3334
$ws.Range("K3").Value = 1954.3334
$ws.Range("M3").Value = -1840.3334

$ws.Range("H94").Value = 5551.625
$ws.Range("I94").Value = 6986.8
$ws.Range("J94").Value = 3159.6667
$ws.Range("K94").Value = 6986.8
$ws.Range("L94").Value = 3159.6667
$ws.Range("M94").Value = -6535.8
$ws.Range("N94").Value = -4061.6667

$ws.Range("H105").Value = 2148.6667
$ws.Range("I105").Value = 2148.5
$ws.Range("J105").Value = 2149
$ws.Range("K105").Value = 2148.5
$ws.Range("L105").Value = 2149
$ws.Range("M105").Value = -401.5
$ws.Range("N105").Value = -5643

$ws.Range("H135").Value = 84849.5
$ws.Range("J135").Value = 84849.5
$ws.Range("L135").Value = 84849.5
$ws.Range("N135").Value = -94989.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("N16").Value = 0
$ws.Range("L16").Value = ""
$ws.Range("M16").Value = ""

$ws.Range("H58").Value = 5497.4
$ws.Range("I58").Value = 2993.5
$ws.Range("K58").Value = 2993.5
$ws.Range("M58").Value = -2790.5

$ws.Range("H62").Value = 12124.25
$ws.Range("I62").Value = 12124.25
$ws.Range("K62").Value = 12124.25
$ws.Range("M62").Value = -11500.25

$ws.Range("H65").Value = 12124.25
$ws.Range("I65").Value = 12124.25
$ws.Range("K65").Value = 60621.25
$ws.Range("M65").Value = -57501.25

$ws.Range("H92").Value = 50000
$ws.Range("J92").Value = 50000
$ws.Range("L92").Value = 50000
$ws.Range("N92").Value = -54992

$ws.Range("H99").Value = 4166.6665
$ws.Range("I99").Value = 3000
$ws.Range("J99").Value = 4400
$ws.Range("K99").Value = 3000
$ws.Range("L99").Value = 4400
$ws.Range("M99").Value = -1502
$ws.Range("N99").Value = -7396

$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("N113").Value = 0
$ws.Range("L113").Value = ""
$ws.Range("M113").Value = ""

$ws.Range("H126").Value = 4166.6665
$ws.Range("I126").Value = 3000
$ws.Range("J126").Value = 4400
$ws.Range("K126").Value = 9000
$ws.Range("L126").Value = 13200
$ws.Range("M126").Value = -6530
$ws.Range("N126").Value = -18140

$ws.Range("H132").Value = 2785.818
$ws.Range("I132").Value = 1992.8572
$ws.Range("J132").Value = 4173.5
$ws.Range("K132").Value = 5978.571599999999
$ws.Range("L132").Value = 12520.5
$ws.Range("M132").Value = -3448.571599999999
$ws.Range("N132").Value = -17580.5

$ws.Range("H136").Value = 5497.4
$ws.Range("I136").Value = 2993.5
$ws.Range("K136").Value = 8980.5
$ws.Range("M136").Value = -6430.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 184
$ws.Range("I38").Value = 144.8
$ws.Range("J38").Value = 249.33333
$ws.Range("K38").Value = 434.4
$ws.Range("L38").Value = 747.99999
$ws.Range("M38").Value = -87.40000000000003
$ws.Range("N38").Value = -1441.99999

$ws.Range("H98").Value = 2830
$ws.Range("J98").Value = 2992.5
$ws.Range("L98").Value = 8977.5
$ws.Range("N98").Value = -11973.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1567.375
$ws.Range("I107").Value = 1448.8
$ws.Range("J107").Value = 1765
$ws.Range("K107").Value = 1448.8
$ws.Range("L107").Value = 1765
$ws.Range("M107").Value = 471.2
$ws.Range("N107").Value = -5605

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3897.8
$ws.Range("I7").Value = 3897.8
$ws.Range("K7").Value = 3897.8
$ws.Range("M7").Value = -3785.8

$ws.Range("H55").Value = 1307
$ws.Range("I55").Value = 634
$ws.Range("J55").Value = 1811.75
$ws.Range("K55").Value = 634
$ws.Range("L55").Value = 1811.75
$ws.Range("M55").Value = -461
$ws.Range("N55").Value = -2157.75

$ws.Range("H100").Value = 2316
$ws.Range("I100").Value = 1861
$ws.Range("J100").Value = 2998.5
$ws.Range("K100").Value = 1861
$ws.Range("L100").Value = 2998.5
$ws.Range("M100").Value = -1320
$ws.Range("N100").Value = -4080.5

$ws.Range("H126").Value = 3897.8
$ws.Range("I126").Value = 3897.8
$ws.Range("K126").Value = 11693.4
$ws.Range("M126").Value = -9223.400000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1182.6
$ws.Range("I122").Value = 1182.6
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3547.8
$ws.Range("L122").Value = 0
$ws.Range("N122").Value = -1097.8
$ws.Range("M122").Value = ""

$ws.Range("H139").Value = 95750
$ws.Range("J139").Value = 95750
$ws.Range("L139").Value = 95750
$ws.Range("N139").Value = -106030
